$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 ("time_taken") - reuse the same bold/centered/bordered
# style already used by the other header cells (e.g. E1) by copying formats.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "time_taken"

# Timestamps for the data rows (2-27), in the order they appear in the diff.
$timestamps = @(
    "2021-10-05 13:38:50.536547",
    "2021-10-05 13:38:50.536555",
    "2021-10-05 13:38:50.536558",
    "2021-10-05 13:38:50.536560",
    "2021-10-05 13:38:50.536562",
    "2021-10-05 13:38:50.536564",
    "2021-10-05 13:38:50.536566",
    "2021-10-05 13:38:50.536568",
    "2021-10-05 13:38:50.536570",
    "2021-10-05 13:38:50.536572",
    "2021-10-05 13:38:50.536574",
    "2021-10-05 13:38:50.536576",
    "2021-10-05 13:38:50.536577",
    "2021-10-05 13:38:50.536579",
    "2021-10-05 13:38:50.536581",
    "2021-10-05 13:38:50.536583",
    "2021-10-05 13:38:50.536585",
    "2021-10-05 13:38:50.536587",
    "2021-10-05 13:38:50.536589",
    "2021-10-05 13:38:50.536591",
    "2021-10-05 13:38:50.536593",
    "2021-10-05 13:38:50.536595",
    "2021-10-05 13:38:50.536597",
    "2021-10-05 13:38:50.536599",
    "2021-10-05 13:38:50.536601",
    "2021-10-05 13:38:50.536603"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
